# Generate Report for Handback
# Two e2e handback items got re-generated with new GUIDs / timestamps:
#   0fa44a94-ddd6-48e1-88bc-b5de6ca30c95  ->  b2fe5179-7a38-43b7-bcbb-1af8a1d01adc
#   9616c1bc-fe0e-4b85-8053-0e21158a9cd4  ->  ffffd63af122-c2f0-423a-a164-e411309fc001
# plus refreshed "latest xliff generate / handoff / handback" timestamps and
# refreshed xlf correspondence file names.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "0fa44a94-ddd6-48e1-88bc-b5de6ca30c95"
$newGuid1 = "b2fe5179-7a38-43b7-bcbb-1af8a1d01adc"
$oldGuid2 = "9616c1bc-fe0e-4b85-8053-0e21158a9cd4"
$newGuid2 = "ffffd63af122-c2f0-423a-a164-e411309fc001"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"

$newXlfZh1 = "$newGuid1.b22865a2cef30841ce65d1e8ceb116d3f935ceeb.zh-cn.xlf"
$newXlfDe1 = "$newGuid1.b22865a2cef30841ce65d1e8ceb116d3f935ceeb.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A2").Value = $newMd1
$wsOv.Range("B2").Value = "e2e\$newMd1"
$wsOv.Range("G2").Value = "2016-09-04 05:09:21"

$wsOv.Range("A3").Value = $newMd2
$wsOv.Range("B3").Value = "e2e\$newMd2"
$wsOv.Range("G3").Value = "2016-09-04 05:09:21"

# refresh hyperlink display text (preserve the original external targets) -
# the runtime collapses the whole-sheet collection on Delete, so wipe once
# and recreate every link on the sheet from scratch, in original order.
$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid1.md", "", "", "e2e\$newMd1") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid2.md", "", "", "e2e\$newMd2") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("I2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh1
$wsZh.Range("H2").Value = "2016-09-04 05:09:16"
$wsZh.Range("J2").Value = $newXlfZh1
$wsZh.Range("K2").Value = "2016-09-04 05:09:34"

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("I3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh1
$wsZh.Range("H3").Value = "2016-09-04 05:09:16"
$wsZh.Range("J3").Value = $newXlfZh1
$wsZh.Range("K3").Value = "2016-09-04 05:09:34"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid1.md", "", "", $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/386ab5d7b266b2b142e850f92c7e2677d6d09d91/e2e/$oldGuid1.md", "", "", $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid2.md", "", "", $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/386ab5d7b266b2b142e850f92c7e2677d6d09d91/e2e/$oldGuid2.md", "", "", $newMd2) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("I2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe1
$wsDe.Range("H2").Value = "2016-09-04 05:09:21"
$wsDe.Range("J2").Value = $newXlfDe1
$wsDe.Range("K2").Value = "2016-09-04 05:09:41"

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("I3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe1
$wsDe.Range("H3").Value = "2016-09-04 05:09:21"
$wsDe.Range("J3").Value = $newXlfDe1
$wsDe.Range("K3").Value = "2016-09-04 05:09:41"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid1.md", "", "", $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a0f180cb97f124065ed96c9c3cdc314d9c40de86/e2e/$oldGuid1.md", "", "", $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/847c88ff7362393cdff2e05ae1cf26330c59832b/e2e/$oldGuid2.md", "", "", $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a0f180cb97f124065ed96c9c3cdc314d9c40de86/e2e/$oldGuid2.md", "", "", $newMd2) | Out-Null
